# Update header/info block
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "Testing"
$ws.Range("C4").Value = "AHU-Horizontal"
$ws.Range("C5").Value = 3000
$ws.Range("D6").Value = "DATE                  :    02-Nov-2020"
$ws.Range("C7").Value = "AHU-510"

# Row 10 - Casing Inner Sheet (part code changes, qty/total recomputed)
$ws.Range("B10").Value = "LCCR07502"
$ws.Range("D10").Value = "GI sheet 0.8mm THK 120 GSM"
$ws.Range("F10").Value = 14.657808
$ws.Range("H10").Value = 43.973423999999994

# Row 11 - Casing Outer Sheet (part code now shared with row10, spec changes)
$ws.Range("B11").Value = "LCCR07502"
$ws.Range("D11").Value = "GI sheet 0.8mm THK 120 GSM"
$ws.Range("F11").Value = 14.657808
$ws.Range("H11").Value = 43.973423999999994

# Row 12 - Corner Profile
$ws.Range("B12").Value = "LAHC23702"
$ws.Range("D12").Value = "ALUMINUM OMEGA PROFILE - 40-23 MM ROUNDED NTBK"
$ws.Range("F12").Value = 3
$ws.Range("H12").Value = 9

# Row 13 - Omega Profile
$ws.Range("B13").Value = "LAHC23714"
$ws.Range("D13").Value = "ALUMINUM CORNER PROFILE - 40 - 45 CHAMFERED TBK"
$ws.Range("F13").Value = 2.7
$ws.Range("H13").Value = 8.100000000000001

# Row 14 - new: Corner Joiner
$ws.Range("A14").Value = 5
$ws.Range("B14").Value = "LAHC15201"
$ws.Range("C14").Value = "Corner Joiner"
$ws.Range("D14").Value = "40mm NTBK Rounded - 3 Way Corner"
$ws.Range("F14").Value = 8
$ws.Range("G14").Value = "Nos"
$ws.Range("H14").Value = 24

# Row 15 - new: Omega Joiner
$ws.Range("A15").Value = 6
$ws.Range("B15").Value = "LAHC15202"
$ws.Range("C15").Value = "Omega Joiner"
$ws.Range("D15").Value = "40mm NTBK Rounded - Omega Jointer"
$ws.Range("F15").Value = 44
$ws.Range("G15").Value = "Nos"
$ws.Range("H15").Value = 132

# Row 16 - new: Polyol
$ws.Range("A16").Value = 7
$ws.Range("B16").Value = "LAHN60001"
$ws.Range("C16").Value = "Polyol"
$ws.Range("D16").Value = "POLYOL"
$ws.Range("F16").Value = 2.585088
$ws.Range("G16").Value = "Kgs"
$ws.Range("H16").Value = 7.7552639999999995

# Row 17 - new: Isol
$ws.Range("A17").Value = 8
$ws.Range("B17").Value = "LAHN60006"
$ws.Range("C17").Value = "Isol"
$ws.Range("D17").Value = "ISOL"
$ws.Range("F17").Value = 1.7233920000000003
$ws.Range("G17").Value = "Kgs"
$ws.Range("H17").Value = 5.1701760000000005
